$d = $word.ActiveDocument
$d.Content.Find.Execute("En esta sección se procederá", $true, $false, $false, $false, $false,
                         $true, 1, $false, "En esta sección, se procederá", 2)
